$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")
$ws.Activate()

# Update sprint backlog status ("in progress" -> "done") for several rows,
# and fill in "Remaining Ressources" (column K) values for a new patient /
# several stories that were missing it.

$ws.Range("K15").Value = 6
$ws.Range("K16").Value = 6
$ws.Range("K17").Value = 2
$ws.Range("K18").Value = 2
$ws.Range("K19").Value = 8
$ws.Range("K21").Value = 4
$ws.Range("K23").Value = 2
$ws.Range("K24").Value = 2

$ws.Range("L14").Value = "done"
$ws.Range("L15").Value = "done"
$ws.Range("L16").Value = "done"
$ws.Range("L17").Value = "done"
$ws.Range("L18").Value = "done"
$ws.Range("L19").Value = "done"
$ws.Range("L21").Value = "done"

# Update the view's scroll position / selection to match where the author
# was working.
$ws.Range("L25").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
